$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wilscher's absence note: "Wilscher krank" -> "Wilscher krank!"
$ws.Range("C16").Value = "Wilscher krank!"

# Sunny's absence note: "Sunny krank" -> "Sunny krank!" (bold, red)
# and ";Menü verbessert" -> "; Menü verbessert" (regular)
$r = $ws.Range("C17")
$r.Value = "Sunny krank!; Menü verbessert"
$c1 = $r.Characters(1, 12)
$c1.Font.Bold = $true
$c1.Font.Color = 255

# Move the active selection to C17 (matches the sheetView selection in the saved file)
$ws.Range("C17").Select()
